$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 11): Prediction / Error updated
$ws.Range("D2").Value = 0.0175080232964091
$ws.Range("E2").Value = 0.0175080232964091

# Row 3 (Control 3): Success flips to TRUE; Prediction / Error updated
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = 0.4961939058923746
$ws.Range("E3").Value = 0.4961939058923746

# Row 4 (Control 38): Prediction / Error updated (very small values)
$ws.Range("D4").Value = [double]"6.588574986123243E-20"
$ws.Range("E4").Value = [double]"6.588574986123243E-20"

# Row 5 (Control 29): Prediction / Error updated
$ws.Range("D5").Value = 0.9980899070973831
$ws.Range("E5").Value = 0.9980899070973831

# Row 6 (MDD 43): Prediction / Error updated
$ws.Range("D6").Value = 0.9929272482003404
$ws.Range("E6").Value = 0.9929272482003404

# Row 7 (MDD 3): Prediction / Error updated
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0

# Row 8 (MDD 19): Prediction / Error updated
$ws.Range("D8").Value = 0.1140219934002382
$ws.Range("E8").Value = 0.8859780065997618

# Row 9 (MDD 7): Prediction / Error updated
$ws.Range("D9").Value = 0.9998026330750298
$ws.Range("E9").Value = 0.0001973669249701882

# Row 11 (MDD 1): Prediction / Error / Cross Entropy Loss / Success % updated
$ws.Range("D11").Value = 0.2925247891852543
$ws.Range("E11").Value = 0.7074752108147457
$ws.Range("F11").Value = 1.531610131263733
$ws.Range("G11").Value = 0.6
